$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.697.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.962.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.52"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.06%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +4.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0795"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.18%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("E12").Value = "  +6.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.833"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.252.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.963.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.643.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "230.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("E26").Value = "  +8.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +17.31%  "
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("E32").Value = "  +4.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0618"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.86%  "
$ws.Range("E35").Value = "  +13.42%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.86%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0982"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0211"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.371.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.90%  "
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.141.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.03%  "
